# Add 3 new groups of data (hyy, hzj, cxq) as columns E, F, G on Sheet1.
# Commit message: "add 3 groups of data hyy hzj cxq"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1) ---
# Written in this order so the new shared-string table entries land in the
# same order as the target workbook (cxq, hyy, hzj).
$ws.Range("E1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("G1").Value = "hyy-调节6Hz_20170306_110203_ASIC_EEG"
$ws.Range("F1").Value = "hzj-调节6Hz_20170220_113105_ASIC_EEG"

# --- New data row 2 ---
$ws.Range("E2").Value = 0.76557863501483681
$ws.Range("F2").Value = 0.81681681681681684
$ws.Range("G2").Value = 0.78787878787878785

# --- New data row 3 ---
$ws.Range("E3").Value = 0.75801749271137031
$ws.Range("F3").Value = 0.76451612903225807
$ws.Range("G3").Value = 0.78498293515358364

# Match the resulting selection state in the saved file: whole column F
# selected with F1 as the active cell.
$ws.Range("F1:F1048576").Select()
